$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 20717.268
$ws.Range("I62").Value = 17704.908
$ws.Range("J62").Value = 29001.25
$ws.Range("K62").Value = 17704.908
$ws.Range("L62").Value = 29001.25
$ws.Range("M62").Value = -17080.908
$ws.Range("N62").Value = -30249.25
$ws.Range("H64").Value = 6914.0454
$ws.Range("I64").Value = 4492.9
$ws.Range("J64").Value = 8931.666999999999
$ws.Range("K64").Value = 4492.9
$ws.Range("L64").Value = 8931.666999999999
$ws.Range("M64").Value = -4244.9
$ws.Range("N64").Value = -9427.666999999999
$ws.Range("H65").Value = 20717.268
$ws.Range("I65").Value = 17704.908
$ws.Range("J65").Value = 29001.25
$ws.Range("K65").Value = 88524.53999999999
$ws.Range("L65").Value = 145006.25
$ws.Range("M65").Value = -85404.53999999999
$ws.Range("N65").Value = -151246.25
$ws.Range("H67").Value = 6914.0454
$ws.Range("I67").Value = 4492.9
$ws.Range("J67").Value = 8931.666999999999
$ws.Range("K67").Value = 4492.9
$ws.Range("L67").Value = 8931.666999999999
$ws.Range("M67").Value = -3634.9
$ws.Range("N67").Value = -10647.667
$ws.Range("H125").Value = 2347.5
$ws.Range("J125").Value = 2347.5
$ws.Range("L125").Value = 21127.5
$ws.Range("N125").Value = -26047.5
$ws.Range("H138").Value = 1747.2632
$ws.Range("I138").Value = 1071.2307
$ws.Range("J138").Value = 3212
$ws.Range("K138").Value = 3213.6921
$ws.Range("L138").Value = 9636
$ws.Range("M138").Value = 1926.3079
$ws.Range("N138").Value = -19916
$ws.Range("H141").Value = 1702.1333
$ws.Range("J141").Value = 2037.5
$ws.Range("L141").Value = 6112.5
$ws.Range("N141").Value = -16472.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6063611
$ws.Range("I32").Value = 6413412
$ws.Range("K32").Value = 6413412
$ws.Range("M32").Value = -6413125
$ws.Range("H61").Value = 8334070
$ws.Range("J61").Value = 999
$ws.Range("L61").Value = 999
$ws.Range("N61").Value = -1423
$ws.Range("H74").Value = 3212.8928
$ws.Range("I74").Value = 2049.762
$ws.Range("J74").Value = 6702.2856
$ws.Range("K74").Value = 2049.762
$ws.Range("L74").Value = 6702.2856
$ws.Range("M74").Value = -1175.762
$ws.Range("N74").Value = -8450.285599999999
$ws.Range("H77").Value = 3212.8928
$ws.Range("I77").Value = 2049.762
$ws.Range("J77").Value = 6702.2856
$ws.Range("K77").Value = 10248.81
$ws.Range("L77").Value = 33511.428
$ws.Range("M77").Value = -5880.810000000001
$ws.Range("N77").Value = -42247.428
$ws.Range("H122").Value = 1703.8334
$ws.Range("I122").Value = 1541.3334
$ws.Range("J122").Value = 1866.3334
$ws.Range("K122").Value = 4624.0002
$ws.Range("L122").Value = 5599.0002
$ws.Range("M122").Value = -2174.0002
$ws.Range("N122").Value = -10499.0002
$ws.Range("H136").Value = 8334070
$ws.Range("J136").Value = 999
$ws.Range("L136").Value = 2997
$ws.Range("N136").Value = -8097
$ws.Range("H139").Value = 92572
$ws.Range("J139").Value = 92572
$ws.Range("L139").Value = 92572
$ws.Range("N139").Value = -102852
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 18034.6
$ws.Range("I44").Value = 15045
$ws.Range("J44").Value = 29993
$ws.Range("K44").Value = 15045
$ws.Range("L44").Value = 29993
$ws.Range("M44").Value = -14548
$ws.Range("N44").Value = -30987
$ws.Range("H134").Value = 1579781.5
$ws.Range("I134").Value = 1703760.2
$ws.Range("K134").Value = 5111280.6
$ws.Range("M134").Value = -5108745.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 63000
$ws.Range("J20").Value = 63000
$ws.Range("L20").Value = 63000
$ws.Range("N20").Value = -63472
$ws.Range("H30").Value = 63000
$ws.Range("J30").Value = 63000
$ws.Range("L30").Value = 63000
$ws.Range("N30").Value = -63182
$ws.Range("H31").Value = 110082.75
$ws.Range("I31").Value = 143045.81
$ws.Range("J31").Value = 37564
$ws.Range("K31").Value = 143045.81
$ws.Range("L31").Value = 37564
$ws.Range("M31").Value = -142750.81
$ws.Range("N31").Value = -38154
$ws.Range("H34").Value = 110082.75
$ws.Range("I34").Value = 143045.81
$ws.Range("J34").Value = 37564
$ws.Range("K34").Value = 143045.81
$ws.Range("L34").Value = 37564
$ws.Range("M34").Value = -142843.81
$ws.Range("N34").Value = -37968
$ws.Range("H58").Value = 4116111.8
$ws.Range("I58").Value = 4116111.8
$ws.Range("K58").Value = 4116111.8
$ws.Range("M58").Value = -4115908.8
$ws.Range("H99").Value = 4345.75
$ws.Range("J99").Value = 4345.75
$ws.Range("L99").Value = 4345.75
$ws.Range("N99").Value = -7341.75
$ws.Range("H126").Value = 4345.75
$ws.Range("J126").Value = 4345.75
$ws.Range("L126").Value = 13037.25
$ws.Range("N126").Value = -17977.25
$ws.Range("H128").Value = 63000
$ws.Range("J128").Value = 63000
$ws.Range("L128").Value = 63000
$ws.Range("N128").Value = -72960
$ws.Range("H132").Value = 216285.05
$ws.Range("J132").Value = 1287030.8
$ws.Range("L132").Value = 3861092.4
$ws.Range("N132").Value = -3866152.4
$ws.Range("H134").Value = 6607.6206
$ws.Range("I134").Value = 6911.4814
$ws.Range("J134").Value = 2505.5
$ws.Range("K134").Value = 20734.4442
$ws.Range("L134").Value = 7516.5
$ws.Range("M134").Value = -18199.4442
$ws.Range("N134").Value = -12586.5
$ws.Range("H136").Value = 4116111.8
$ws.Range("I136").Value = 4116111.8
$ws.Range("K136").Value = 12348335.4
$ws.Range("M136").Value = -12345785.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4590.8
$ws.Range("I70").Value = 4620.8125
$ws.Range("J70").Value = 4470.75
$ws.Range("K70").Value = 4620.8125
$ws.Range("L70").Value = 4470.75
$ws.Range("M70").Value = -4350.8125
$ws.Range("N70").Value = -5010.75
$ws.Range("H73").Value = 4590.8
$ws.Range("I73").Value = 4620.8125
$ws.Range("J73").Value = 4470.75
$ws.Range("K73").Value = 4620.8125
$ws.Range("L73").Value = 4470.75
$ws.Range("M73").Value = -3684.8125
$ws.Range("N73").Value = -6342.75
$ws.Range("H80").Value = 304839.25
$ws.Range("I80").Value = 369089.66
$ws.Range("J80").Value = 5004
$ws.Range("K80").Value = 369089.66
$ws.Range("L80").Value = 5004
$ws.Range("M80").Value = -368091.66
$ws.Range("N80").Value = -7000
$ws.Range("H83").Value = 304839.25
$ws.Range("I83").Value = 369089.66
$ws.Range("J83").Value = 5004
$ws.Range("K83").Value = 1845448.3
$ws.Range("L83").Value = 25020
$ws.Range("M83").Value = -1840456.3
$ws.Range("N83").Value = -35004
$ws.Range("H93").Value = 80223
$ws.Range("J93").Value = 80223
$ws.Range("L93").Value = 80223
$ws.Range("N93").Value = -83967
$ws.Range("H102").Value = 3082.0952
$ws.Range("I102").Value = 2266.6897
$ws.Range("J102").Value = 4901.077
$ws.Range("K102").Value = 2266.6897
$ws.Range("L102").Value = 4901.077
$ws.Range("M102").Value = -644.6896999999999
$ws.Range("N102").Value = -8145.077
$ws.Range("H122").Value = 8358
$ws.Range("I122").Value = 4373.7
$ws.Range("J122").Value = 14998.5
$ws.Range("K122").Value = 13121.1
$ws.Range("L122").Value = 44995.5
$ws.Range("M122").Value = -10671.1
$ws.Range("N122").Value = -49895.5
$ws.Range("H132").Value = 718198.8
$ws.Range("I132").Value = 1013066.3
$ws.Range("K132").Value = 3039198.9
$ws.Range("M132").Value = -3036668.9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4050.5293
$ws.Range("I7").Value = 3918.5
$ws.Range("K7").Value = 3918.5
$ws.Range("M7").Value = -3806.5
$ws.Range("H40").Value = 4570.8335
$ws.Range("I40").Value = 4545.5884
$ws.Range("K40").Value = 4545.5884
$ws.Range("M40").Value = -4409.5884
$ws.Range("H68").Value = 3858.4285
$ws.Range("I68").Value = 4500
$ws.Range("K68").Value = 4500
$ws.Range("M68").Value = -3751
$ws.Range("H71").Value = 3858.4285
$ws.Range("I71").Value = 4500
$ws.Range("K71").Value = 22500
$ws.Range("M71").Value = -18756
$ws.Range("H82").Value = 197.33333
$ws.Range("I82").Value = 197.33333
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 197.33333
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 163.66667
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 197.33333
$ws.Range("I85").Value = 197.33333
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 197.33333
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 1050.66667
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 4050.5293
$ws.Range("I126").Value = 3918.5
$ws.Range("K126").Value = 11755.5
$ws.Range("M126").Value = -9285.5
$ws.Range("H130").Value = 89959.42999999999
$ws.Range("J130").Value = 89959.42999999999
$ws.Range("L130").Value = 89959.42999999999
$ws.Range("N130").Value = -99999.42999999999
$ws.Range("H136").Value = 105992.914
$ws.Range("I136").Value = 2370
$ws.Range("J136").Value = 251065
$ws.Range("K136").Value = 7110
$ws.Range("L136").Value = 753195
$ws.Range("M136").Value = -4560
$ws.Range("N136").Value = -758295
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3800
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3800
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5048
$ws.Range("H65").Value = 3800
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 19000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -25240
$ws.Range("H107").Value = 1683.7037
$ws.Range("I107").Value = 939.94446
$ws.Range("J107").Value = 3171.2222
$ws.Range("K107").Value = 2819.83338
$ws.Range("L107").Value = 9513.6666
$ws.Range("M107").Value = -899.83338
$ws.Range("N107").Value = -13353.6666
$ws.Range("H126").Value = 6998.1816
$ws.Range("I126").Value = 6997.1665
$ws.Range("K126").Value = 20991.4995
$ws.Range("M126").Value = -18521.4995
$ws.Range("H132").Value = 5442694.5
$ws.Range("I132").Value = 5921814.5
$ws.Range("J132").Value = 12666.667
$ws.Range("K132").Value = 17765443.5
$ws.Range("L132").Value = 38000.001
$ws.Range("M132").Value = -17762913.5
$ws.Range("N132").Value = -43060.001
